# Add Peter's file info to the "Features" sheet, row 2 (the "user winning
# teams selection" task), mirroring the existing entries in that column
# (e.g. row 4/5/6: Who / Which files are touched).

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Features")
$ws.Activate()

$ws.Range("C2").Value = "Peter"
$ws.Range("D2").Value = "MarchMadnessGUI, Bracket"

# Update the selected cell shown in the sheet view to C8, matching the
# saved cursor position after making the edit.
[void]$ws.Range("C8").Select()
